# Auto-generated edit script applying the Gilgamesh_Profits market-data refresh
# (currentAveragePrice* / LevePrice* / LeveProfit* columns H..N) per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1947.7333
$ws.Range("I43").Value = 1696.5
$ws.Range("K43").Value = 1696.5
$ws.Range("M43").Value = -1627.5
$ws.Range("H55").Value = 475.55554
$ws.Range("I55").Value = 763.6667
$ws.Range("K55").Value = 763.6667
$ws.Range("M55").Value = -549.6667
$ws.Range("H103").Value = 4465621
$ws.Range("J103").Value = 11905206
$ws.Range("L103").Value = 35715618
$ws.Range("N103").Value = -35716790
$ws.Range("H116").Value = 4187.375
$ws.Range("I116").Value = 4153.231
$ws.Range("K116").Value = 4153.231
$ws.Range("M116").Value = -711.2309999999998
$ws.Range("H132").Value = 6471.2915
$ws.Range("I132").Value = 6709.174
$ws.Range("K132").Value = 20127.522
$ws.Range("M132").Value = -17597.522
$ws.Range("H137").Value = 2181745.8
$ws.Range("I137").Value = 2943563.8
$ws.Range("K137").Value = 8830691.399999999
$ws.Range("M137").Value = -8828141.399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 112960.32
$ws.Range("I74").Value = 130830.28
$ws.Range("J74").Value = 3187.7144
$ws.Range("K74").Value = 130830.28
$ws.Range("L74").Value = 3187.7144
$ws.Range("M74").Value = -129956.28
$ws.Range("N74").Value = -4935.7144
$ws.Range("H77").Value = 112960.32
$ws.Range("I77").Value = 130830.28
$ws.Range("J77").Value = 3187.7144
$ws.Range("K77").Value = 654151.4
$ws.Range("L77").Value = 15938.572
$ws.Range("M77").Value = -649783.4
$ws.Range("N77").Value = -24674.572
$ws.Range("H110").Value = 9508.75
$ws.Range("I110").Value = 9340.666999999999
$ws.Range("K110").Value = 9340.666999999999
$ws.Range("M110").Value = -7295.666999999999
$ws.Range("H132").Value = 1999.5769
$ws.Range("I132").Value = 1349.0769
$ws.Range("K132").Value = 4047.2307
$ws.Range("M132").Value = -1517.2307

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 252.54546
$ws.Range("J80").Value = 255.33333
$ws.Range("L80").Value = 255.33333
$ws.Range("N80").Value = -2251.33333
$ws.Range("H83").Value = 252.54546
$ws.Range("J83").Value = 255.33333
$ws.Range("L83").Value = 1276.66665
$ws.Range("N83").Value = -11260.66665
$ws.Range("H86").Value = 3558.0476
$ws.Range("I86").Value = 3224.6
$ws.Range("J86").Value = 4391.6665
$ws.Range("K86").Value = 3224.6
$ws.Range("L86").Value = 4391.6665
$ws.Range("M86").Value = -2101.6
$ws.Range("N86").Value = -6637.6665
$ws.Range("H89").Value = 3558.0476
$ws.Range("I89").Value = 3224.6
$ws.Range("J89").Value = 4391.6665
$ws.Range("K89").Value = 16123
$ws.Range("L89").Value = 21958.3325
$ws.Range("M89").Value = -10507
$ws.Range("N89").Value = -33190.3325
$ws.Range("H107").Value = 3078377.5
$ws.Range("I107").Value = 4274918.5
$ws.Range("J107").Value = 1557.8572
$ws.Range("K107").Value = 4274918.5
$ws.Range("L107").Value = 1557.8572
$ws.Range("M107").Value = -4272998.5
$ws.Range("N107").Value = -5397.8572
$ws.Range("H134").Value = 2192.4866
$ws.Range("I134").Value = 1530.7391
$ws.Range("J134").Value = 3279.6428
$ws.Range("K134").Value = 4592.2173
$ws.Range("L134").Value = 9838.928400000001
$ws.Range("M134").Value = -2057.2173
$ws.Range("N134").Value = -14908.9284

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1999
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H31").Value = 4755.109
$ws.Range("I31").Value = 3317.3845
$ws.Range("K31").Value = 3317.3845
$ws.Range("M31").Value = -3022.3845
$ws.Range("H34").Value = 4755.109
$ws.Range("I34").Value = 3317.3845
$ws.Range("K34").Value = 3317.3845
$ws.Range("M34").Value = -3115.3845
$ws.Range("H58").Value = 2215.1333
$ws.Range("I58").Value = 1476
$ws.Range("K58").Value = 1476
$ws.Range("M58").Value = -1273
$ws.Range("H99").Value = 3344
$ws.Range("I99").Value = 2682.6667
$ws.Range("K99").Value = 2682.6667
$ws.Range("M99").Value = -1184.6667
$ws.Range("H113").Value = 1999
$ws.Range("J113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("N113").ClearContents()
$ws.Range("H126").Value = 3344
$ws.Range("I126").Value = 2682.6667
$ws.Range("K126").Value = 8048.000100000001
$ws.Range("M126").Value = -5578.000100000001
$ws.Range("H132").Value = 15158136
$ws.Range("I132").Value = 6444.5454
$ws.Range("K132").Value = 19333.6362
$ws.Range("M132").Value = -16803.6362
$ws.Range("H134").Value = 7947.3
$ws.Range("I134").Value = 8163.6665
$ws.Range("K134").Value = 24490.9995
$ws.Range("M134").Value = -21955.9995
$ws.Range("H136").Value = 2215.1333
$ws.Range("I136").Value = 1476
$ws.Range("K136").Value = 4428
$ws.Range("M136").Value = -1878

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 989.7037
$ws.Range("J2").Value = 1483.7059
$ws.Range("L2").Value = 8902.2354
$ws.Range("N2").Value = -9128.2354
$ws.Range("H38").Value = 341.66666
$ws.Range("J38").Value = 515.5714
$ws.Range("L38").Value = 1546.7142
$ws.Range("N38").Value = -2240.7142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2506.5
$ws.Range("J113").Value = 2506.5
$ws.Range("L113").Value = 2506.5
$ws.Range("N113").Value = -6846.5
$ws.Range("H132").Value = 2308
$ws.Range("J132").Value = 2957.1428
$ws.Range("L132").Value = 8871.428400000001
$ws.Range("N132").Value = -13931.4284

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 11267.143
$ws.Range("I40").Value = 11267.143
$ws.Range("K40").Value = 11267.143
$ws.Range("M40").Value = -11131.143
$ws.Range("H55").Value = 612.2917
$ws.Range("I55").Value = 491.64285
$ws.Range("J55").Value = 781.2
$ws.Range("K55").Value = 491.64285
$ws.Range("L55").Value = 781.2
$ws.Range("M55").Value = -318.64285
$ws.Range("N55").Value = -1127.2
$ws.Range("H93").Value = 450.55554
$ws.Range("I93").Value = 422.85715
$ws.Range("J93").Value = 547.5
$ws.Range("K93").Value = 422.85715
$ws.Range("L93").Value = 547.5
$ws.Range("M93").Value = 825.14285
$ws.Range("N93").Value = -3043.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 17758
$ws.Range("J41").Value = 17307.6
$ws.Range("L41").Value = 17307.6
$ws.Range("N41").Value = -18087.6
$ws.Range("H107").Value = 788.7083
$ws.Range("I107").Value = 878.9474
$ws.Range("J107").Value = 445.8
$ws.Range("K107").Value = 2636.8422
$ws.Range("L107").Value = 1337.4
$ws.Range("M107").Value = -716.8422
$ws.Range("N107").Value = -5177.4
$ws.Range("H122").Value = 7355883
$ws.Range("I122").Value = 3083.9
$ws.Range("K122").Value = 9251.700000000001
$ws.Range("M122").Value = -6801.700000000001
